$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 - flight to Stockholm (ARN)
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "Sunday, Jan 08"
$ws.Range("C46").Value = "8:10 PM"
$ws.Range("D46").Value = "FR1944"
$ws.Range("E46").Value = "Stockholm"
$ws.Range("F46").Value = "(ARN)"
$ws.Range("G46").Value = "Ryanair "
$ws.Range("H46").Value = "B738"
$ws.Range("I46").Value = "(SP-RKL)"
$ws.Range("J46").Value = "8:37 PM"
$ws.Range("K46").Font.Size = 11
$ws.Range("L46").Value = "0 hours, 27 minutes"
$ws.Range("M46").Font.Size = 11

# Row 47 - flight to Vienna (VIE)
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "Sunday, Jan 08"
$ws.Range("C47").Value = "8:15 PM"
$ws.Range("D47").Value = "FR1574"
$ws.Range("E47").Value = "Vienna"
$ws.Range("F47").Value = "(VIE)"
$ws.Range("G47").Value = "Ryanair "
$ws.Range("H47").Value = "B738"
$ws.Range("I47").Value = "(SP-RKT)"
$ws.Range("J47").Value = "8:20 PM"
$ws.Range("K47").Font.Size = 11
$ws.Range("L47").Value = "0 hours, 5 minutes"
$ws.Range("M47").Font.Size = 11

# Row 48 - flight to London (STN)
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "Sunday, Jan 08"
$ws.Range("C48").Value = "9:40 PM"
$ws.Range("D48").Value = "FR2670"
$ws.Range("E48").Value = "London"
$ws.Range("F48").Value = "(STN)"
$ws.Range("G48").Value = "Ryanair "
$ws.Range("H48").Value = "B738"
$ws.Range("I48").Value = "(EI-EKS)"
$ws.Range("J48").Value = "9:52 PM"
$ws.Range("K48").Font.Size = 11
$ws.Range("L48").Value = "0 hours, 12 minutes"
$ws.Range("M48").Font.Size = 11
